$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('LP1912')
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws3 = $wb.Worksheets.Item('6203-6173')

# ---- LP1912 ----
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 06:58:01'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 60'
$ws1.Cells.Item(18, 1).Value = '04:56:49'
$ws1.Cells.Item(18, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(18, 4).Value = 68
$ws1.Cells.Item(19, 1).Value = '05:23:04'
$ws1.Cells.Item(19, 3).Value = '10_OLMOS'
$ws1.Cells.Item(19, 4).Value = 41
$ws1.Cells.Item(37, 1).Value = '06:58:01'
$ws1.Cells.Item(37, 4).Value = 3
$ws1.Cells.Item(38, 1).Value = '06:58:01'
$ws1.Cells.Item(38, 4).Value = 6
$ws1.Cells.Item(40, 1).Value = '06:58:01'
$ws1.Cells.Item(40, 4).Value = 9
$ws1.Cells.Item(42, 1).Value = '06:58:01'
$ws1.Cells.Item(42, 4).Value = 16
$ws1.Cells.Item(43, 1).Value = '06:58:01'
$ws1.Cells.Item(43, 4).Value = 23
$ws1.Cells.Item(44, 1).Value = '06:58:01'
$ws1.Cells.Item(44, 4).Value = 26
$ws1.Cells.Item(45, 1).Value = '06:58:01'
$ws1.Cells.Item(45, 4).Value = 31
$ws1.Cells.Item(47, 1).Value = '06:58:01'
$ws1.Cells.Item(47, 4).Value = 36
$ws1.Cells.Item(49, 1).Value = '06:58:01'
$ws1.Cells.Item(49, 4).Value = 38
$ws1.Cells.Item(50, 1).Value = '06:58:01'
$ws1.Cells.Item(50, 4).Value = 39
$ws1.Cells.Item(52, 1).Value = '06:58:01'
$ws1.Cells.Item(52, 4).Value = 46
$ws1.Cells.Item(53, 1).Value = '06:58:01'
$ws1.Cells.Item(53, 4).Value = 51
$ws1.Cells.Item(56, 1).Value = '06:58:01'
$ws1.Cells.Item(56, 4).Value = 62
$ws1.Cells.Item(57, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(58, 1).Value = '06:58:01'
$ws1.Cells.Item(58, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(58, 4).Value = 65
$ws1.Cells.Item(59, 1).Value = '06:58:01'
$ws1.Cells.Item(59, 2).Value = '08:04'
$ws1.Cells.Item(59, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(59, 4).Value = 66
$ws1.Cells.Item(60, 1).Value = '06:58:01'
$ws1.Cells.Item(60, 2).Value = '08:14'
$ws1.Cells.Item(60, 3).Value = '10_OLMOS'
$ws1.Cells.Item(60, 4).Value = 76
$ws1.Cells.Item(61, 1).Value = '06:58:01'
$ws1.Cells.Item(61, 2).Value = '08:19'
$ws1.Cells.Item(61, 3).Value = '15_ABASTO'
$ws1.Cells.Item(61, 4).Value = 81
$ws1.Cells.Item(62, 1).Value = '06:58:01'
$ws1.Cells.Item(62, 2).Value = '08:30'
$ws1.Cells.Item(62, 3).Value = '14_ABASTO'
$ws1.Cells.Item(62, 4).Value = 92
$ws1.Cells.Item(62, 5).Value = 'LP1912'
$ws1.Cells.Item(63, 1).Value = '06:58:01'
$ws1.Cells.Item(63, 2).Value = '08:34'
$ws1.Cells.Item(63, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(63, 4).Value = 96
$ws1.Cells.Item(63, 5).Value = 'LP1912'
$ws1.Cells.Item(64, 1).Value = '06:58:01'
$ws1.Cells.Item(64, 2).Value = '08:48'
$ws1.Cells.Item(64, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(64, 4).Value = 110
$ws1.Cells.Item(64, 5).Value = 'LP1912'
$ws1.Cells.Item(65, 1).Value = '06:58:01'
$ws1.Cells.Item(65, 2).Value = '08:51'
$ws1.Cells.Item(65, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(65, 4).Value = 113
$ws1.Cells.Item(65, 5).Value = 'LP1912'

# ---- LP1912-215 ----
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 06:58:01'
$ws2.Cells.Item(3, 1).Value = 'Total filas: 12'
$ws2.Cells.Item(14, 1).Value = '06:58:01'
$ws2.Cells.Item(14, 4).Value = 9
$ws2.Cells.Item(15, 1).Value = '06:58:01'
$ws2.Cells.Item(15, 4).Value = 23
$ws2.Cells.Item(16, 1).Value = '06:58:01'
$ws2.Cells.Item(16, 4).Value = 96
$ws2.Cells.Item(17, 1).Value = '06:58:01'
$ws2.Cells.Item(17, 2).Value = '08:48'
$ws2.Cells.Item(17, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(17, 4).Value = 110
$ws2.Cells.Item(17, 5).Value = 'LP1912'

# ---- 6203-6173 ----
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 06:58:01'
$ws3.Cells.Item(3, 1).Value = 'Total filas: 9'
$ws3.Cells.Item(11, 1).Value = '06:58:01'
$ws3.Cells.Item(11, 2).Value = '07:36'
$ws3.Cells.Item(11, 4).Value = 38
$ws3.Cells.Item(12, 2).Value = '08:10'
$ws3.Cells.Item(12, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(12, 4).Value = 84
$ws3.Cells.Item(12, 5).Value = 'L6173'
$ws3.Cells.Item(13, 1).Value = '06:58:01'
$ws3.Cells.Item(13, 2).Value = '08:23'
$ws3.Cells.Item(13, 3).Value = '215C_LA PLATA'
$ws3.Cells.Item(13, 4).Value = 85
$ws3.Cells.Item(13, 5).Value = 'L6203'
$ws3.Cells.Item(14, 1).Value = '06:58:01'
$ws3.Cells.Item(14, 2).Value = '08:52'
$ws3.Cells.Item(14, 3).Value = '215A_LA PLATA'
$ws3.Cells.Item(14, 4).Value = 114
$ws3.Cells.Item(14, 5).Value = 'L6173'
